# Auto-generated edit script: updates odds values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.83
$ws.Range("H2").Value = 3.5
$ws.Range("I2").Value = 4.75
$ws.Range("L2").Value = 5
$ws.Range("Z2").Value = 8
$ws.Range("AG2").Value = 17
$ws.Range("AI2").Value = 351
$ws.Range("AJ2").Value = 12
$ws.Range("AK2").Value = 23
$ws.Range("AM2").Value = 51

# Row 3
$ws.Range("G3").Value = 1.22
$ws.Range("H3").Value = 7.5
$ws.Range("I3").Value = 10
$ws.Range("J3").Value = 1.57
$ws.Range("K3").Value = 3.1
$ws.Range("L3").Value = 8.5
$ws.Range("N3").Value = 26
$ws.Range("Q3").Value = 1.3
$ws.Range("R3").Value = 3.5
$ws.Range("S3").Value = 1.8
$ws.Range("T3").Value = 2
$ws.Range("W3").Value = 1.7
$ws.Range("X3").Value = 2.05
$ws.Range("Z3").Value = 9
$ws.Range("AB3").Value = 9
$ws.Range("AD3").Value = 21
$ws.Range("AF3").Value = 15
$ws.Range("AG3").Value = 21
$ws.Range("AI3").Value = 151
$ws.Range("AL3").Value = 26
$ws.Range("AN3").Value = 51
$ws.Range("AO3").Value = 41

# Row 4
$ws.Range("G4").Value = 1.8
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 2.6
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("S4").Value = 5.5
$ws.Range("T4").Value = 1.14
$ws.Range("W4").Value = 2.5
$ws.Range("X4").Value = 1.5
$ws.Range("Y4").Value = 4.75
$ws.Range("AB4").Value = 13
$ws.Range("AC4").Value = 19
$ws.Range("AJ4").Value = 9.5
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 19
$ws.Range("AO4").Value = 67
$ws.Range("AP4").Value = 4.6
$ws.Range("AQ4").Value = 1.2
$ws.Range("AR4").Value = 2.05
$ws.Range("AS4").Value = 1.8

# Row 5
$ws.Range("J5").Value = 3.75
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 5.5
$ws.Range("U5").Value = 1.67
$ws.Range("V5").Value = 2.1

# Row 17
$ws.Range("O17").Value = 1.5
$ws.Range("P17").Value = 2.5
$ws.Range("Q17").Value = 2.6
$ws.Range("R17").Value = 1.48
$ws.Range("AR17").Value = 2
$ws.Range("AS17").Value = 1.85

# Row 19
$ws.Range("G19").Value = 3
$ws.Range("I19").Value = 2.4
$ws.Range("J19").Value = 3.5
$ws.Range("K19").Value = 2.1
$ws.Range("M19").Value = 1.06
$ws.Range("N19").Value = 10
$ws.Range("O19").Value = 1.29
$ws.Range("P19").Value = 3.5
$ws.Range("Q19").Value = 1.98
$ws.Range("R19").Value = 1.88
$ws.Range("S19").Value = 3.25
$ws.Range("T19").Value = 1.33
$ws.Range("W19").Value = 1.7
$ws.Range("X19").Value = 2.05
$ws.Range("AC19").Value = 23
$ws.Range("AE19").Value = 10
$ws.Range("AF19").Value = 6
$ws.Range("AG19").Value = 13
$ws.Range("AI19").Value = 201
$ws.Range("AJ19").Value = 8.5
$ws.Range("AK19").Value = 12
$ws.Range("AL19").Value = 9.5
$ws.Range("AO19").Value = 29

# Row 20
$ws.Range("M20").Value = 1.03
$ws.Range("N20").Value = 15
$ws.Range("O20").Value = 1.2
$ws.Range("P20").Value = 4.33

# Row 24
$ws.Range("G24").Value = 1.55
$ws.Range("I24").Value = 5.75
$ws.Range("W24").Value = 1.91
$ws.Range("X24").Value = 1.8
$ws.Range("AH24").Value = 51
$ws.Range("AI24").Value = 351

# Row 28
$ws.Range("O28").Value = 1.18
$ws.Range("P28").Value = 4.5
$ws.Range("Q28").Value = 1.6
$ws.Range("R28").Value = 2.3
$ws.Range("S28").Value = 2.5
$ws.Range("T28").Value = 1.5
$ws.Range("AP28").Value = 2
$ws.Range("AQ28").Value = 1.85

# Row 29
$ws.Range("G29").Value = 1.5
$ws.Range("I29").Value = 5
$ws.Range("J29").Value = 1.95
$ws.Range("N29").Value = 26
$ws.Range("Q29").Value = 1.3
$ws.Range("R29").Value = 3.5
$ws.Range("U29").Value = 1.18
$ws.Range("V29").Value = 4.5
$ws.Range("AB29").Value = 13
$ws.Range("AH29").Value = 34
$ws.Range("AJ29").Value = 26
$ws.Range("AO29").Value = 29

# Row 32
$ws.Range("G32").Value = 1.95
$ws.Range("I32").Value = 3.75
$ws.Range("J32").Value = 2.6
$ws.Range("L32").Value = 4
$ws.Range("M32").Value = 1.04
$ws.Range("N32").Value = 13
$ws.Range("Q32").Value = 1.83
$ws.Range("R32").Value = 2.03

# Row 33
$ws.Range("G33").Value = 4.75
$ws.Range("H33").Value = 3.8
$ws.Range("I33").Value = 1.7
$ws.Range("J33").Value = 5.5
$ws.Range("K33").Value = 2.2
$ws.Range("L33").Value = 2.3
$ws.Range("W33").Value = 2
$ws.Range("X33").Value = 1.73
$ws.Range("Z33").Value = 23
$ws.Range("AF33").Value = 7.5
$ws.Range("AG33").Value = 19
$ws.Range("AH33").Value = 67
$ws.Range("AK33").Value = 7.5

# Row 34
$ws.Range("G34").Value = 2.2
$ws.Range("I34").Value = 3.5
$ws.Range("J34").Value = 3.1
$ws.Range("AC34").Value = 21
$ws.Range("AK34").Value = 17
$ws.Range("AR34").Value = 1.95
$ws.Range("AS34").Value = 1.9

# Row 39
$ws.Range("G39").Value = 2.05
$ws.Range("I39").Value = 3.8
$ws.Range("N39").Value = 10
$ws.Range("O39").Value = 1.29
$ws.Range("P39").Value = 3.5
$ws.Range("Q39").Value = 2
$ws.Range("R39").Value = 1.85
$ws.Range("Z39").Value = 9.5
$ws.Range("AB39").Value = 17
$ws.Range("AD39").Value = 26
$ws.Range("AG39").Value = 15
$ws.Range("AI39").Value = 251
$ws.Range("AO39").Value = 41

# Row 41
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 4
$ws.Range("I41").Value = 1.75
$ws.Range("J41").Value = 4
$ws.Range("L41").Value = 2.3
$ws.Range("Q41").Value = 1.5
$ws.Range("R41").Value = 2.5
$ws.Range("S41").Value = 2.2
$ws.Range("T41").Value = 1.62
$ws.Range("Y41").Value = 17
$ws.Range("Z41").Value = 23
$ws.Range("AA41").Value = 13
$ws.Range("AF41").Value = 8
$ws.Range("AG41").Value = 12
$ws.Range("AN41").Value = 13

# Row 42
$ws.Range("J42").Value = 3.1
$ws.Range("O42").Value = 1.25
$ws.Range("P42").Value = 3.75
$ws.Range("S42").Value = 3
$ws.Range("T42").Value = 1.36
$ws.Range("AG42").Value = 13
$ws.Range("AJ42").Value = 10
$ws.Range("AO42").Value = 29

# Row 43
$ws.Range("G43").Value = 2.6
$ws.Range("I43").Value = 2.6
$ws.Range("J43").Value = 3.25
